$wb = $excel.ActiveWorkbook
$survey = $wb.Worksheets.Item("survey")
$promptTypes = $wb.Worksheets.Item("prompt_types")

# --- survey sheet: insert a new row 3 for the phone_number text prompt ---
# (this pushes the existing send_sms row from row 3 down to row 4)
$survey.Rows.Item(3).Insert()

# Fill the new row 3. Typing order matters for the shared-string table, so
# we create the new unique strings in the exact order they must land in
# xl/sharedStrings.xml: odk_sms, the sms-bridge text, phone_number, the
# enter-phone-number text, odk_sms_automatic, the no-confirmation text.
$survey.Range("A5").Value = "odk_sms"
$survey.Range("C5").Value = "This will send an sms via the sms bridge."

$survey.Range("A3").Value = "text"
$survey.Range("B3").Value = "phone_number"
$survey.Range("C3").Value = "Enter the phone number to which to send the text."

$survey.Range("A6").Value = "odk_sms_automatic"
$survey.Range("C6").Value = "This will send an sms without requiring confirmation."

$survey.Range("B5").Value = "odk_sms"
$survey.Range("B6").Value = "odk_sms_automatic"

# Row heights for the new rows.
$survey.Rows.Item(3).RowHeight = 12
$survey.Rows.Item(5).RowHeight = 12.75
$survey.Rows.Item(6).RowHeight = 12.75

# Column B needs to widen to fit "phone_number".
$survey.Columns.Item(2).ColumnWidth = 15.25

# --- prompt_types sheet: register the two new prompt types ---
$promptTypes.Range("A3").Value = "odk_sms"
$promptTypes.Range("B3").Value = "integer"
$promptTypes.Range("A4").Value = "odk_sms_automatic"
$promptTypes.Range("B4").Value = "integer"

# --- selections / active sheet ---
$survey.Range("C7").Select()
$promptTypes.Activate()
$promptTypes.Range("B5").Select()
